# adding sputtering prevac attributes
# - Insert new "ALD" and "ALD_BeneQ" sheets (schema stub sheets) right after
#   "OzoneCleaning" and before "CVD".
# - Rename "sputtering" -> "Sputtering" and "sputtering_prevac" -> "Sputtering_prevac".
# - Expand the "Sputtering_prevac" header row with the full set of process
#   attribute columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the two new stub sheets ("ALD", "ALD_BeneQ") after "OzoneCleaning"
#    and before "CVD", matching the simple name/id/iri schema used by the
#    other process sheets.
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("OzoneCleaning")

$ald = $wb.Worksheets.Add($null, $afterSheet)
$ald.Name = "ALD"
$ald.Range("A1").Value = "name"
$ald.Range("B1").Value = "id"
$ald.Range("C1").Value = "iri"

$aldBeneQ = $wb.Worksheets.Add($null, $ald)
$aldBeneQ.Name = "ALD_BeneQ"
$aldBeneQ.Range("A1").Value = "name"
$aldBeneQ.Range("B1").Value = "id"
$aldBeneQ.Range("C1").Value = "iri"

# ---------------------------------------------------------------------------
# 2. Rename the sputtering sheets to their capitalized forms.
# ---------------------------------------------------------------------------
$sputtering = $wb.Worksheets.Item("sputtering")
$sputtering.Name = "Sputtering"

$sputteringPrevac = $wb.Worksheets.Item("sputtering_prevac")
$sputteringPrevac.Name = "Sputtering_prevac"

# ---------------------------------------------------------------------------
# 3. Replace the "Sputtering_prevac" header row with the full attribute list.
# ---------------------------------------------------------------------------
$columns = @(
    "sputtering_prevac_id",
    "substrate_id",
    "sample_owner",
    "process_user",
    "date",
    "holder",
    "notes",
    "step_number",
    "orientation",
    "sputter_pressure",
    "substrate_temperature",
    "ramp",
    "rotation",
    "z_position",
    "gas",
    "flow_rate",
    "target_position",
    "target",
    "target_power",
    "DC_RF",
    "time",
    "name",
    "id",
    "iri"
)

for ($i = 0; $i -lt $columns.Length; $i++) {
    $sputteringPrevac.Cells.Item(1, $i + 1).Value = $columns[$i]
}
